# Repull data, push all data, mean calculation
# Update the dSF (column F) values for the rows whose dSF figure changed
# after re-pulling / recalculating the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8  = 2
    10 = -1
    11 = -3
    12 = -3
    13 = -2
    14 = 3
    18 = 2
    27 = 0
    28 = 1
    33 = 0
    34 = -5
    37 = -3
    38 = 2
    40 = 2
    42 = 0
    52 = -2
    57 = -3
    58 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
